$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The phone-number column ("no_hp_mitra", currently G) carries a special
# (hyperlink-like) cell format. Two new columns are being inserted before
# it, so that format needs to move two columns to the right (to I) while G
# itself reverts to the plain style used by its neighbours.
$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

# New header labels (status_pekerjaan / detail_pekerjaan) take over columns
# G and H, pushing no_hp_mitra / email_mitra / tahun one column further
# right (G->I, H->J, I->K).
$ws.Range("G1").Value = "status_pekerjaan"
$ws.Range("H1").Value = "detail_pekerjaan"
$ws.Range("I1").Value = "no_hp_mitra"
$ws.Range("J1").Value = "email_mitra"
$ws.Range("K1").Value = "tahun"
$ws.Range("K1").NumberFormat = "@"

# Sample row values shifted the same way.
$ws.Range("G2").Value = "1"
$ws.Range("H2").Value = "contoh"
$ws.Range("I2").Value = "+62"
$ws.Range("J2").Value = "contoh@gmail.com"
$ws.Range("K2").Value = "16-02-2024"
$ws.Range("K2").NumberFormat = "@"

# sobat_id sample value changed.
$ws.Range("A2").Value = "1122"
